# Updated cryptos list (Price / Volume(1h) columns) to reflect the latest
# scrape. Price values (column D) are entered with a leading apostrophe so
# Excel keeps them as literal text (e.g. "599.96", "0.0000109") instead of
# silently re-parsing them as numbers and dropping formatting / precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.523.29"
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = "'2.658.00"
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'599.96"
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").Value = "'156.78"
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +5.12%  '
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").Value = "'5.86"
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").Value = "'29.25"
$ws.Range("E14").Value = '  -4.97%  '
$ws.Range("D15").Value = "'3.134.63"
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").Value = "'65.325.94"
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").Value = "'2.662.12"
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = "'12.63"
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("D20").Value = "'7.58"
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").Value = "'350.39"
$ws.Range("E21").Value = '  -2.53%  '
$ws.Range("D23").Value = "'69.25"
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").Value = "'0.0000109"
$ws.Range("E24").Value = '  +3.45%  '
$ws.Range("D25").Value = "'9.69"
$ws.Range("E25").Value = '  +1.28%  '
$ws.Range("E26").Value = '  -5.14%  '
$ws.Range("D27").Value = "'0.167"
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("E28").Value = '  -3.50%  '
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = "'534.00"
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = "'2.13"
$ws.Range("E32").Value = '  -4.82%  '
$ws.Range("E33").Value = '  -2.10%  '
$ws.Range("D34").Value = "'6.50"
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("D35").Value = "'5.46"
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("D36").Value = "'0.422"
$ws.Range("E36").Value = '  -2.54%  '
$ws.Range("D37").Value = "'20.35"
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("D39").Value = "'159.21"
$ws.Range("E40").Value = '  -3.49%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").Value = "'42.66"
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("D43").Value = "'164.20"
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("E44").Value = '  -2.27%  '
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("E47").Value = '  -1.02%  '
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("E49").Value = '  -2.32%  '
$ws.Range("E50").Value = '  +3.02%  '
$ws.Range("D51").Value = "'20.11"
$ws.Range("E51").Value = '  +1.53%  '
